# Update Name of Algo
# Apply updated RandomForest imputed values to column D for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -7.384999999999995
    13 = -8.374900000000002
    16 = -8.787300000000005
    18 = -8.589599999999997
    20 = -7.234100000000004
    26 = -8.679700000000002
    27 = -8.785499999999999
    29 = -7.317799999999998
    35 = -7.661700000000005
    36 = -7.760900000000003
    45 = -7.5642
    55 = -8.9598
    57 = -8.461600000000004
    69 = -7.089299999999993
    76 = -7.4263
    78 = -7.590300000000004
    82 = -8.287699999999994
    83 = -9.1225
    93 = -6.376599999999996
    97 = -8.082099999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
